$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ntf5"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3590443333333333
$ws.Range("H2").Value = 1.077133
$ws.Range("I2").Value = 0.480597973884934
$ws.Range("J2").Value = 0.480597973884934
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.293221
$ws.Range("N2").Value = 0.879663
$ws.Range("O2").Value = 0.02611983441994871
$ws.Range("P2").Value = 0.02611983441994871
$ws.Range("Q2").Value = 0.1052793384643333
$ws.Range("R2").Value = 0.9475140461789999
$ws.Range("S2").Value = 0.01255313950043731
$ws.Range("T2").Value = 0.01255313950043731

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntf5"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3590443333333333
$ws.Range("H3").Value = 1.077133
$ws.Range("I3").Value = 0.480597973884934
$ws.Range("J3").Value = 0.480597973884934
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.74310933333333
$ws.Range("N3").Value = 32.229328
$ws.Range("O3").Value = 0.9569854715114954
$ws.Range("P3").Value = 0.9569854715114955
$ws.Range("Q3").Value = 3.857252528513777
$ws.Range("R3").Value = 34.71527275662399
$ws.Range("S3").Value = 0.4599252786457429
$ws.Range("T3").Value = 0.459925278645743

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntf5"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3590443333333333
$ws.Range("H4").Value = 1.077133
$ws.Range("I4").Value = 0.480597973884934
$ws.Range("J4").Value = 0.480597973884934
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.183427
$ws.Range("N4").Value = 0.550281
$ws.Range("O4").Value = 0.01633949433413
$ws.Range("P4").Value = 0.01633949433413
$ws.Range("Q4").Value = 0.06585842493033332
$ws.Range("R4").Value = 0.5927258243729999
$ws.Range("S4").Value = 0.007852727871287235
$ws.Range("T4").Value = 0.007852727871287237

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf5"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3590443333333333
$ws.Range("H5").Value = 1.077133
$ws.Range("I5").Value = 0.480597973884934
$ws.Range("J5").Value = 0.480597973884934
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006232666666666667
$ws.Range("N5").Value = 0.018698
$ws.Range("O5").Value = 0.0005551997344257983
$ws.Range("P5").Value = 0.0005551997344257983
$ws.Range("Q5").Value = 0.002237803648222222
$ws.Range("R5").Value = 0.020140232834
$ws.Range("S5").Value = 0.0002668278674664921
$ws.Range("T5").Value = 0.0002668278674664921

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Ntf5"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.388034
$ws.Range("H6").Value = 1.164102
$ws.Range("I6").Value = 0.519402026115066
$ws.Range("J6").Value = 0.5194020261150661
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.293221
$ws.Range("N6").Value = 0.879663
$ws.Range("O6").Value = 0.02611983441994871
$ws.Range("P6").Value = 0.02611983441994871
$ws.Range("Q6").Value = 0.113779717514
$ws.Range("R6").Value = 1.024017457626
$ws.Range("S6").Value = 0.0135666949195114
$ws.Range("T6").Value = 0.01356669491951141

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Ntf5"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.388034
$ws.Range("H7").Value = 1.164102
$ws.Range("I7").Value = 0.519402026115066
$ws.Range("J7").Value = 0.5194020261150661
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.74310933333333
$ws.Range("N7").Value = 32.229328
$ws.Range("O7").Value = 0.9569854715114954
$ws.Range("P7").Value = 0.9569854715114955
$ws.Range("Q7").Value = 4.168691687050666
$ws.Range("R7").Value = 37.518225183456
$ws.Range("S7").Value = 0.4970601928657525
$ws.Range("T7").Value = 0.4970601928657526

$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Ntf5"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.388034
$ws.Range("H8").Value = 1.164102
$ws.Range("I8").Value = 0.519402026115066
$ws.Range("J8").Value = 0.5194020261150661
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.183427
$ws.Range("N8").Value = 0.550281
$ws.Range("O8").Value = 0.01633949433413
$ws.Range("P8").Value = 0.01633949433413
$ws.Range("Q8").Value = 0.071175912518
$ws.Range("R8").Value = 0.640583212662
$ws.Range("S8").Value = 0.008486766462842762
$ws.Range("T8").Value = 0.008486766462842766

$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Ntf5"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.388034
$ws.Range("H9").Value = 1.164102
$ws.Range("I9").Value = 0.519402026115066
$ws.Range("J9").Value = 0.5194020261150661
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.006232666666666667
$ws.Range("N9").Value = 0.018698
$ws.Range("O9").Value = 0.0005551997344257983
$ws.Range("P9").Value = 0.0005551997344257983
$ws.Range("Q9").Value = 0.002418486577333333
$ws.Range("R9").Value = 0.021766379196
$ws.Range("S9").Value = 0.0135666949195114
$ws.Range("T9").Value = 0.0002883718669593063

$ws.Rows.Item(10).Delete()